# Insert the new "Directives" section into the document, just before the
# trailing empty paragraph at the very end of the body (which stays last,
# unchanged, exactly as in the target diff).
$d = $word.ActiveDocument

$lastParaIndex = $d.Paragraphs.Count
$insertionRange = $d.Paragraphs($lastParaIndex).Range
$insertionRange.Collapse(1)

$newSectionXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:highlight w:val="yellow"/></w:rPr><w:t>Directives</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> [change behaviours of DOM element]</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="7"/></w:numPr><w:ind w:left="284"/></w:pPr><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>Structural directive</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p/><w:p><w:r><w:t>*</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ngIf</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> - &gt;Add / remove element from the DOM</w:t></w:r></w:p><w:p><w:r><w:t>*</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ngFor</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> -&gt; dynamic DOM structure update</w:t></w:r></w:p><w:p><w:r><w:t>*</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ngSwitch</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> -&gt; Switch on basis of conditions</w:t></w:r></w:p><w:p/><w:p><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve">Angular 18+ </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>updates :</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve"> need to import </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>commonModule</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve"> separately to use directives</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$null = $insertionRange.InsertXML($newSectionXml)

Write-Output "Inserted Directives section; paragraph count now $($d.Paragraphs.Count)"
